$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.0637
$ws.Range("A9").Value = -21.81260000000001
$ws.Range("C12").Value = -11.1536
$ws.Range("C14").Value = -13.62039999999999
$ws.Range("A18").Value = -22.31220000000001
$ws.Range("A20").Value = -21.46429999999998
$ws.Range("C26").Value = -12.89150000000001
$ws.Range("A27").Value = -21.84829999999999
$ws.Range("C27").Value = -12.5084
$ws.Range("C29").Value = -11.1962
$ws.Range("A35").Value = -21.03859999999998
$ws.Range("C37").Value = -14.23649999999999
$ws.Range("C38").Value = -13.2285
$ws.Range("C51").Value = -12.4178
$ws.Range("C52").Value = -11.2228
$ws.Range("C55").Value = -13.6163
$ws.Range("A69").Value = -21.8112
$ws.Range("C69").Value = -11.2373
$ws.Range("C70").Value = -11.7637
$ws.Range("A76").Value = -19.88489999999999
$ws.Range("A78").Value = -19.94449999999998
$ws.Range("C81").Value = -13.0419
$ws.Range("A82").Value = -21.9861
$ws.Range("A83").Value = -21.9544
$ws.Range("C83").Value = -12.9888
$ws.Range("A93").Value = -21.225
$ws.Range("C102").Value = -13.639
